$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge the two runs around "and then run ----- " / "in order
# to set up your database. " into a single run (no text change, just a
# formatting-identical run merge, which Word performs naturally when a
# Find/Replace spans the run boundary).
# ---------------------------------------------------------------------
$text1 = "and then run ----- in order to set up your database."
$d.Content.Find.Execute($text1, $false, $false, $false, $false, $false, $true, 1, $false, $text1, 2) | Out-Null

# ---------------------------------------------------------------------
# Change 2: merge "How to configure the con" + "nection string" into a
# single run "How to configure the connection string".
# ---------------------------------------------------------------------
$text2 = "How to configure the connection string"
$d.Content.Find.Execute($text2, $false, $false, $false, $false, $false, $true, 1, $false, $text2, 2) | Out-Null

# ---------------------------------------------------------------------
# Change 3: expand the "In the web.config file, ..." paragraph with
# additional clarifying text (Open ... then ... (this is inbetween
# <connectionStrings> and </connectionStrings>), modify the
# connectionString for SchedContext  (found between <add name=
# "SchedContext" and />) to match your database credentials.)
# ---------------------------------------------------------------------
$old3 = "In the web.config file, in the connection strings section of the file, modify the connectionString for SchedContext to match your database credentials."
$rng = $d.Content
$found3 = $rng.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $start = $rng.Start
    $end = $rng.End
    $target = $d.Range($start, $end)
    $target.Text = ""
    $insertPoint = $d.Range($start, $start)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/></w:rPr><w:t>Open</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/></w:rPr><w:t>web.config</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> file,</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> then</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> in the connection strings section of the file</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> (this is </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/></w:rPr><w:t>inbetween</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>connectionStrings</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>&gt;</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>&lt;</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>connectionStrings</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>&gt;</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> modify the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/></w:rPr><w:t>connectionString</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> for </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>SchedContext</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">found between </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>&lt;</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>add</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="FF0000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>name</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>=</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>"</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="0000FF"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>SchedContext</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>"</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> and /&gt;</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>to match your</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:i/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> database credentials.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertPoint.InsertXML($xml)
}
